$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the last data row (previously row 13 / Accident 91, dated 2018-11-25).
# The sheet now only covers 2018-11-26, one fewer accident record than before.
$ws.Range("A13").EntireRow.Delete()

# Date/Time/Hour columns hold text-like values (e.g. "21", "2018-11-26") that Excel
# would otherwise auto-coerce into dates/numbers, so mark them as Text first.
$dateRange = $ws.Range("F2:F12")
$timeRange = $ws.Range("G2:G12")
$hourRange = $ws.Range("P2:P12")
$dateRange.NumberFormat = "@"
$timeRange.NumberFormat = "@"
$hourRange.NumberFormat = "@"

# Row 2
$ws.Range("A2").Value = 1
$ws.Range("C2").Value = "Injuries"
$ws.Range("D2").Value = 35.244997
$ws.Range("E2").Value = -85.10848799999999
$ws.Range("F2").Value = "2018-11-26"
$ws.Range("G2").Value = "21:29:37"
$ws.Range("H2").Value = "10731-10759 Hixson Pike"
$ws.Range("K2").Value = "HAMILTON COUNTY"
$ws.Range("P2").Value = "21"
$ws.Range("Y2").Value = 11

# Row 3
$ws.Range("A3").Value = 1
$ws.Range("C3").Value = "Injuries"
$ws.Range("D3").Value = 35.244997
$ws.Range("E3").Value = -85.10848799999999
$ws.Range("F3").Value = "2018-11-26"
$ws.Range("G3").Value = "21:29:37"
$ws.Range("H3").Value = "10731-10759 Hixson Pike"
$ws.Range("K3").Value = "HAMILTON COUNTY"
$ws.Range("P3").Value = "21"
$ws.Range("Y3").Value = 11

# Row 4
$ws.Range("A4").Value = 22
$ws.Range("C4").Value = "Injuries"
$ws.Range("D4").Value = 35.098369
$ws.Range("E4").Value = -85.327973
$ws.Range("F4").Value = "2018-11-26"
$ws.Range("G4").Value = "17:16:02"
$ws.Range("H4").Value = "Mountain Creek Rd / Signal Mountain Rd"
$ws.Range("K4").Value = "CHATTANOOGA"
$ws.Range("P4").Value = "17"
$ws.Range("Y4").Value = 11

# Row 5
$ws.Range("A5").Value = 23
$ws.Range("C5").Value = "Injuries"
$ws.Range("D5").Value = 35.098369
$ws.Range("E5").Value = -85.327973
$ws.Range("F5").Value = "2018-11-26"
$ws.Range("G5").Value = "17:15:52"
$ws.Range("H5").Value = "Mountain Creek Rd / Signal Mountain Rd"
$ws.Range("K5").Value = "CHATTANOOGA"
$ws.Range("P5").Value = "17"
$ws.Range("Y5").Value = 11

# Row 6
$ws.Range("A6").Value = 25
$ws.Range("C6").Value = "No Injuries"
$ws.Range("D6").Value = 35.075019
$ws.Range("E6").Value = -85.061744
$ws.Range("F6").Value = "2018-11-26"
$ws.Range("G6").Value = "16:56:50"
$ws.Range("H6").Value = "Main St / Ocoee St"
$ws.Range("K6").Value = "COLLEGEDALE"
$ws.Range("P6").Value = "16"
$ws.Range("Y6").Value = 11

# Row 7
$ws.Range("A7").Value = 28
$ws.Range("C7").Value = "Injuries"
$ws.Range("D7").Value = 35.032718
$ws.Range("E7").Value = -85.263344
$ws.Range("F7").Value = "2018-11-26"
$ws.Range("G7").Value = "16:39:30"
$ws.Range("H7").Value = "Glenwood Dr / Oak St"
$ws.Range("K7").Value = "CHATTANOOGA"
$ws.Range("P7").Value = "16"
$ws.Range("Y7").Value = 11

# Row 8
$ws.Range("A8").Value = 29
$ws.Range("C8").Value = "Injuries"
$ws.Range("D8").Value = 35.032718
$ws.Range("E8").Value = -85.263344
$ws.Range("F8").Value = "2018-11-26"
$ws.Range("G8").Value = "16:38:25"
$ws.Range("H8").Value = "Glenwood Dr / Oak St"
$ws.Range("K8").Value = "CHATTANOOGA"
$ws.Range("P8").Value = "16"
$ws.Range("Y8").Value = 11

# Row 9
$ws.Range("A9").Value = 35
$ws.Range("C9").Value = "Injuries"
$ws.Range("D9").Value = 35.032049
$ws.Range("E9").Value = -85.311255
$ws.Range("F9").Value = "2018-11-26"
$ws.Range("G9").Value = "16:21:36"
$ws.Range("H9").Value = "100 W 20th St"
$ws.Range("K9").Value = "CHATTANOOGA"
$ws.Range("P9").Value = "16"
$ws.Range("Y9").Value = 11

# Row 10
$ws.Range("A10").Value = 71
$ws.Range("C10").Value = "Injuries"
$ws.Range("D10").Value = 35.008246
$ws.Range("E10").Value = -85.20027399999999
$ws.Range("F10").Value = "2018-11-26"
$ws.Range("G10").Value = "07:26:33"
$ws.Range("H10").Value = "200 Interstate 75 Sb"
$ws.Range("K10").Value = "CHATTANOOGA"
$ws.Range("P10").Value = "7"
$ws.Range("Y10").Value = 11

# Row 11
$ws.Range("A11").Value = 72
$ws.Range("C11").Value = "Injuries"
$ws.Range("D11").Value = 35.008246
$ws.Range("E11").Value = -85.20027399999999
$ws.Range("F11").Value = "2018-11-26"
$ws.Range("G11").Value = "07:26:33"
$ws.Range("H11").Value = "200 Interstate 75 Sb"
$ws.Range("K11").Value = "CHATTANOOGA"
$ws.Range("P11").Value = "7"
$ws.Range("Y11").Value = 11

# Row 12
$ws.Range("A12").Value = 85
$ws.Range("C12").Value = "Injuries"
$ws.Range("D12").Value = 35.015821
$ws.Range("E12").Value = -85.144578
$ws.Range("F12").Value = "2018-11-26"
$ws.Range("G12").Value = "00:08:13"
$ws.Range("H12").Value = "1715 JENKINS RD"
$ws.Range("K12").Value = "CHATTANOOGA"
$ws.Range("P12").Value = "0"
$ws.Range("Y12").Value = 11

# Drop the temporary Text format now that the values are stored as strings,
# so the cells end up unstyled just like the rest of the data rows.
$dateRange.ClearFormats()
$timeRange.ClearFormats()
$hourRange.ClearFormats()
